# seuk_04.xlsx - "updated data and Kubas photo"
# (this workbook only contains the spreadsheet-data portion of that commit;
#  the photo lives in a different file in the original repo)
#
# Adds the next round-robin results (rows 60-75, match date 2023-04-18 /
# serial 45034) below the existing data, plus a couple of small formatting
# side effects that came along with the resave (row 46/47 height tweak,
# a stray formatted-but-empty O59 cell, and moving the active selection).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- new match rows: Date, Attack_1, Defence_1, Attack_2, Defence_2, G1, G2 ---
$rows = @(
  @(45034, "Krzysiu", "Dawid",  "Kuba",   "Wojtek", 3, 8),
  @(45034, "Dawid",   "Kuba",   "Szymon", "Krzysiu",3, 8),
  @(45034, "Kuba",    "Szymon", "Krzysiu","Wojtek", 8, 2),
  @(45034, "Kuba",    "Dawid",  "Wojtek", "Szymon", 4, 8),
  @(45034, "Szymon",  "Eryk",   "Dawid",  "Krzysiu",6, 8),
  @(45034, "Eryk",    "Kuba",   "Wojtek", "Dawid",  8, 2),
  @(45034, "Krzysiu", "Eryk",   "Szymon", "Kuba",   2, 8),
  @(45034, "Wojtek",  "Krzysiu","Dawid",  "Szymon", 7, 8),
  @(45034, "Kuba",    "Krzysiu","Eryk",   "Szymon", 5, 8),
  @(45034, "Szymon",  "Wojtek", "Dawid",  "Eryk",   8, 4),
  @(45034, "Wojtek",  "Szymon", "Eryk",   "Dawid",  8, 3),
  @(45034, "Wojtek",  "Dawid",  "Eryk",   "Szymon", 8, 4),
  @(45034, "Szymon",  "Dawid",  "Wojtek", "Eryk",   8, 2),
  @(45034, "Eryk",    "Szymon", "Dawid",  "Wojtek", 8, 4),
  @(45034, "Wojtek",  "Szymon", "Eryk",   "Dawid",  3, 8),
  @(45034, "Szymon",  "Wojtek", "Dawid",  "Eryk",   8, 7)
)

$r = 60
foreach ($row in $rows) {
    for ($c = 1; $c -le 7; $c++) {
        $ws.Cells.Item($r, $c).Value = $row[$c - 1]
    }
    $ws.Cells.Item($r, 8).Formula = "=IF(F$r>G$r,1,2)"
    $r = $r + 1
}

# stray formatted (but empty) cell at O59 - same date number format as column A
$ws.Range("O59").NumberFormat = "[$-415]YYYY\-MM\-DD"

# rows 46/47 picked up the 13.8pt "compact" row height used by every row from 48 on
$ws.Rows.Item(46).RowHeight = 13.8
$ws.Rows.Item(47).RowHeight = 13.8

# the brand-new rows (63-75) use that same 13.8pt height too (60-62 already existed
# as empty placeholder rows at 13.8pt, so they keep it automatically)
for ($rr = 63; $rr -le 75; $rr++) {
    $ws.Rows.Item($rr).RowHeight = 13.8
}

# move the active selection the way the author left it
$ws.Range("K70").Select()
